# "WF Real. arranca nuevamente. v006 a v010"
# Mark rows 4 and 6 (A4=3, A6=5) as already run ("Sí" in column H,
# the "FE.Corrido" column) -- matching rows 2 and 3 which are already
# marked, then leave the selection on I1 (just past the used range).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("H4").Value = "Sí"
$ws.Range("H6").Value = "Sí"

$ws.Range("I1").Select()
